$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 48
$prev = $row - 1

$ws.Cells.Item($row, 1).Value2 = "YES"
$ws.Cells.Item($row, 2).Value2 = 48
$ws.Cells.Item($row, 3).Value2 = "KrrSimple"
$ws.Cells.Item($row, 4).Value2 = "Nimsoft"
$ws.Cells.Item($row, 5).Value2 = "Nim01"
$ws.Cells.Item($row, 6).Value2 = ":::URL::www.magenta.ca|order"
$ws.Cells.Item($row, 7).Value2 = "NoAnswer"
$ws.Cells.Item($row, 8).Value2 = "MAJOR"
$ws.Cells.Item($row, 9).Value2 = "'"
$ws.Cells.Item($row, 10).Value2 = "'"
$ws.Cells.Item($row, 11).Value2 = 1386604685
$ws.Cells.Item($row, 12).Value2 = 1386604685
$ws.Cells.Item($row, 13).Value2 = 0
$ws.Cells.Item($row, 14).Value2 = 1386604685

# Copy formats from the row above so the new row matches existing styling
# (the B column uses a bold quote-prefixed style, I/J a quote-prefixed style).
$ws.Cells.Item($prev, 2).Copy()
$ws.Cells.Item($row, 2).PasteSpecial(-4122)

$ws.Cells.Item($prev, 9).Copy()
$ws.Cells.Item($row, 9).PasteSpecial(-4122)

$ws.Cells.Item($prev, 10).Copy()
$ws.Cells.Item($row, 10).PasteSpecial(-4122)

$ws.Application.CutCopyMode = $false

# Leave the selection just below the newly added row, matching where a user
# would land after typing the new entry.
[void]$ws.Range("B50").Select()
